$d = $word.ActiveDocument
$full = $d.Content
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p w:rsidR="00EF1F9E" w:rsidRDefault="003C73FC">
<w:r><w:t>September 7, 2016:  Setup of initial project, using individual user accounts for login and authentication (gold level); we’ve included it in the initial setup since Visual Studio will generate this functionality for us.</w:t></w:r>
<w:r><w:t xml:space="preserve">  Created a separate file for “priority” </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>enum</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>.</w:t></w:r>
</w:p>
<w:p w:rsidR="00437724" w:rsidRDefault="00437724">
<w:r><w:t>Also, the bronze feature of “clearing items” refers to the shopping list, not individual items.</w:t></w:r>
</w:p>
<w:p w:rsidR="00437724" w:rsidRDefault="00437724">
<w:r><w:t>Note from yesterday:  “color” property is to change the color of the text of the list and its contents.</w:t></w:r>
</w:p>
<w:p w:rsidR="00437724" w:rsidRDefault="00437724"/>
<w:p w:rsidR="00437724" w:rsidRDefault="00437724">
<w:r><w:t xml:space="preserve">September 8, 2016:  added the Note property to the </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ShoppingListItem</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> Model as a string, instead of setting it up as its own class.  We feel that this will make it more convenient for the end user.</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p w:rsidR="00437724" w:rsidRDefault="00437724"/>
<w:sectPr w:rsidR="00437724"><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@
$full.InsertXML($xml)
Write-Output "Paragraphs: $($d.Paragraphs.Count)"
